$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Rename the worksheet tab from the generic "1" to "Borjomi"
# ------------------------------------------------------------------
$ws.Name = "Borjomi"

# ------------------------------------------------------------------
# 2. The "Urban" (row 6) and "Rural" (row 7) breakdown for 2010-2023
#    is no longer published for this table - every year (columns
#    B:O) now shows the confidential/unavailable placeholder that is
#    already used by several cells ("..." ellipsis marker). The
#    "Total" row (row 5) is left untouched.
# ------------------------------------------------------------------
$placeholder = [char]0x2026

$urbanRow = 6
$ruralRow = 7
for ($col = 2; $col -le 15; $col++) {
    $ws.Cells.Item($urbanRow, $col).Value = $placeholder
    $ws.Cells.Item($ruralRow, $col).Value = $placeholder
}

# ------------------------------------------------------------------
# 3. Normalize the placeholder glyph everywhere on the sheet: the
#    single-character ellipsis "…" becomes three literal dots "..."
# ------------------------------------------------------------------
[void]$ws.Cells.Replace($placeholder, "...")

# ------------------------------------------------------------------
# 4. Remove the stray blank row that used to separate the data table
#    from the footnote, so the note moves from row 9 up to row 8.
# ------------------------------------------------------------------
$ws.Rows(8).Delete()
